$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.781.06'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.492.84'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '587.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.515'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.141'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.77%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.18%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.340'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.73%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.949.08'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.69'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.15%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '67.699.61'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.09%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000171'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.84%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.511.24'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.24%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.53'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.62%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.96'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.23%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '350.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.07%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.04'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.27'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.70%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.75'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -2.17%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.12'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.24%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.615.15'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.982'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.82%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0902'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.37%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '505.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.76%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.82'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.26'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.63%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.77'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.57%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.00'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.01%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.122'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +4.15%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '162.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.84%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.32'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.51%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.34'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.42%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.00'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.74'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.80%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.328'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.37%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '4.85'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.96%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.41'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.67%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '144.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.52'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.515'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.12%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0254'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0743'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.10%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.59'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.586'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.61%  '
